$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 2 (shifts existing row 2 down to row 4)
$ws.Rows.Item(2).Resize(2).Insert()

# The inserted rows pick up the header row's style; clear formatting so the
# new data rows end up with the default (unstyled) appearance, matching the
# rest of the data rows.
$ws.Rows.Item(2).Resize(2).ClearFormats()

# Force the new cells to be stored as text (not numbers), matching the rest
# of the sheet where numeric-looking scores are kept as strings.
$ws.Range("A2:M3").NumberFormat = "@"

# Row 2: Acessar_agência_virtual.txt
$ws.Range("A2").Value = "Acessar_agência_virtual.txt"
$ws.Range("B2").Value = "5"
$ws.Range("C2").Value = "5"
$ws.Range("D2").Value = "2"
$ws.Range("E2").Value = "5"
$ws.Range("F2").Value = "5"
$ws.Range("G2").Value = "5"
$ws.Range("H2").Value = "4"
$ws.Range("I2").Value = "2"
$ws.Range("J2").Value = "1.5"
$ws.Range("K2").Value = "1.5"
$ws.Range("L2").Value = "2"
$ws.Range("M2").Value = "0"

# Row 3: Acessar_dados_do_portal_da_transparência.txt
$ws.Range("A3").Value = "Acessar_dados_do_portal_da_transparência.txt"
$ws.Range("B3").Value = "5"
$ws.Range("C3").Value = "5"
$ws.Range("D3").Value = "2"
$ws.Range("E3").Value = "5"
$ws.Range("F3").Value = "5"
$ws.Range("G3").Value = "5"
$ws.Range("H3").Value = "4"
$ws.Range("I3").Value = "2"
$ws.Range("J3").Value = "1.5"
$ws.Range("K3").Value = "1.5"
$ws.Range("L3").Value = "2"
$ws.Range("M3").Value = "0"

# Drop the explicit "Text" number format now that the values are locked in as
# text, so the new cells end up with the default (unstyled) appearance
# instead of carrying a leftover custom style index.
$ws.Range("A2:M3").ClearFormats()

$wb.Save()
